$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 1.67
$ws.Range("O2").Value = 2.1
$ws.Range("G3").Value = 4.1
$ws.Range("H3").Value = 4.1
$ws.Range("I3").Value = 1.75
$ws.Range("N3").Value = 1.47
$ws.Range("U3").Value = 26
$ws.Range("V3").Value = 15
$ws.Range("Y3").Value = 29
$ws.Range("AH3").Value = 15
$ws.Range("AI3").Value = 12
$ws.Range("N4").Value = 2.15
$ws.Range("O4").Value = 1.62
$ws.Range("N6").Value = 2.01
$ws.Range("O6").Value = 1.77
$ws.Range("I8").Value = 1.91
$ws.Range("N8").Value = 1.92
$ws.Range("O8").Value = 1.82
$ws.Range("J9").Value = 1.05
$ws.Range("K9").Value = 11
$ws.Range("N9").Value = 1.92
$ws.Range("O9").Value = 1.82
$ws.Range("J10").Value = 1.05
$ws.Range("K10").Value = 11
$ws.Range("N10").Value = 1.79
$ws.Range("O10").Value = 1.94
$ws.Range("Z10").Value = 11
$ws.Range("AD10").Value = 201
$ws.Range("G11").Value = 3.05
$ws.Range("I11").Value = 2.27
$ws.Range("M11").Value = 2.37
$ws.Range("N11").Value = 2.32
$ws.Range("O11").Value = 1.47
$ws.Range("R11").Value = 2
$ws.Range("S11").Value = 1.65
$ws.Range("T11").Value = 7.3
$ws.Range("U11").Value = 14
$ws.Range("V11").Value = 11.75
$ws.Range("W11").Value = 37
$ws.Range("X11").Value = 32
$ws.Range("Y11").Value = 50
$ws.Range("AC11").Value = 110
$ws.Range("AE11").Value = 6.1
$ws.Range("AF11").Value = 9.75
$ws.Range("AG11").Value = 9.75
$ws.Range("AH11").Value = 22
$ws.Range("AI11").Value = 22
$ws.Range("G13").Value = 4.33
$ws.Range("I13").Value = 1.8
$ws.Range("J13").Value = 1.08
$ws.Range("K13").Value = 8
$ws.Range("L13").Value = 1.44
$ws.Range("M13").Value = 2.63
$ws.Range("N13").Value = 2.35
$ws.Range("O13").Value = 1.57
$ws.Range("P13").Value = 1.5
$ws.Range("Q13").Value = 2.5
$ws.Range("R13").Value = 2.1
$ws.Range("S13").Value = 1.67
$ws.Range("T13").Value = 9.5
$ws.Range("V13").Value = 15
$ws.Range("W13").Value = 51
$ws.Range("Y13").Value = 51
$ws.Range("Z13").Value = 8
$ws.Range("AB13").Value = 21
$ws.Range("AD13").Value = 501
$ws.Range("AE13").Value = 5.5
$ws.Range("AF13").Value = 7.5
$ws.Range("G14").Value = 8.25
$ws.Range("H14").Value = 3.9
$ws.Range("L14").Value = 1.33
$ws.Range("M14").Value = 3.05
$ws.Range("N14").Value = 1.93
$ws.Range("O14").Value = 1.7
$ws.Range("R14").Value = 2.37
$ws.Range("S14").Value = 1.52
$ws.Range("T14").Value = 14
$ws.Range("V14").Value = 22
$ws.Range("W14").Value = 175
$ws.Range("X14").Value = 90
$ws.Range("Y14").Value = 80
$ws.Range("Z14").Value = 8.25
$ws.Range("AB14").Value = 19.5
$ws.Range("AC14").Value = 110
$ws.Range("AE14").Value = 4.6
$ws.Range("AF14").Value = 4.8
$ws.Range("AG14").Value = 7.3
$ws.Range("AH14").Value = 7
$ws.Range("AI14").Value = 10.5
$ws.Range("AJ14").Value = 29
$ws.Range("P16").Value = 1.41
$ws.Range("Q16").Value = 2.7
$ws.Range("R16").Value = 1.88
$ws.Range("S16").Value = 1.83
$ws.Range("N20").Value = 1.43
$ws.Range("O20").Value = 2.65
$ws.Range("R20").Value = 2.11
$ws.Range("S20").Value = 1.65
$ws.Range("I22").Value = 2.6
$ws.Range("T22").Value = 10.25
$ws.Range("U22").Value = 14
$ws.Range("X22").Value = 18
$ws.Range("AE22").Value = 10
$ws.Range("AF22").Value = 14
$ws.Range("AJ22").Value = 26
$ws.Range("W24").Value = 41
$ws.Range("J25").Value = 1.04
$ws.Range("L25").Value = 1.33
$ws.Range("P25").Value = 1.44
$ws.Range("Q25").Value = 2.63
$ws.Range("J26").Value = 1.04
$ws.Range("K26").Value = 8
$ws.Range("L26").Value = 1.27
$ws.Range("N26").Value = 2
$ws.Range("O26").Value = 1.8
$ws.Range("P26").Value = 1.4
$ws.Range("J29").ClearContents()
$ws.Range("K29").ClearContents()
$ws.Range("L29").Value = 1.03
$ws.Range("N29").Value = 1.17
$ws.Range("N30").Value = 1.75
$ws.Range("N32").Value = 1.8
$ws.Range("O33").Value = 1.5
